$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 111 (REFERENCIA / G111 is filled in later, after row 118 - see below)
$ws.Range("A111").Value = 3
$ws.Range("B111").Value = 90036484392
$ws.Range("C111").Value = "CIMPRE"
$ws.Range("D111").Value = "SALUD OCUPACIONAL S.A.S."
$ws.Range("E111").Value = 51
$ws.Range("F111").Value = "CC"
$ws.Range("H111").Value = 87451

# Row 112
$ws.Range("A112").Value = 3
$ws.Range("B112").Value = 8903043450
$ws.Range("C112").Value = "ELECTRICOS"
$ws.Range("D112").Value = "DEL VALLE SA"
$ws.Range("E112").Value = 51
$ws.Range("F112").Value = "CC"
$ws.Range("G112").Value = "4565465465"
$ws.Range("H112").Value = 172500

# Row 113
$ws.Range("A113").Value = 3
$ws.Range("B113").Value = 830076882
$ws.Range("C113").Value = "Hp"
$ws.Range("D113").Value = "Financial Services Colombia LLC Sucursal Colombia"
$ws.Range("E113").Value = 7
$ws.Range("F113").Value = "CA"
$ws.Range("G113").Value = "4654654"
$ws.Range("H113").Value = 4407849

# Row 114
$ws.Range("A114").Value = 3
$ws.Range("B114").Value = 1143940722
$ws.Range("C114").Value = "IZC"
$ws.Range("D114").Value = "MAYORISTA SAS"
$ws.Range("E114").Value = 7
$ws.Range("F114").Value = "CA"
$ws.Range("G114").Value = "4565654"
$ws.Range("H114").Value = 52092009

# Row 115
$ws.Range("A115").Value = 3
$ws.Range("B115").Value = 900892841
$ws.Range("C115").Value = "LILIUM"
$ws.Range("D115").Value = "TECNOLOGIA SAS"
$ws.Range("E115").Value = 7
$ws.Range("F115").Value = "CC"
$ws.Range("G115").Value = "465654"
$ws.Range("H115").Value = 669600

# Row 116
$ws.Range("A116").Value = 3
$ws.Range("B116").Value = 800035776
$ws.Range("C116").Value = "NEXSYS"
$ws.Range("D116").Value = "DE COLOMBIA SA"
$ws.Range("E116").Value = 7
$ws.Range("F116").Value = "CA"
$ws.Range("G116").Value = "6546546"
$ws.Range("H116").Value = 18089916

# Row 117
$ws.Range("A117").Value = 3
$ws.Range("B117").Value = 830034343
$ws.Range("C117").Value = "RENTEK"
$ws.Range("D117").Value = "SAS"
$ws.Range("E117").Value = 7
$ws.Range("F117").Value = "CC"
$ws.Range("G117").Value = "465464"
$ws.Range("H117").Value = 4094318

# Row 118
$ws.Range("A118").Value = 3
$ws.Range("B118").Value = 800179308
$ws.Range("C118").Value = "YAMAKI"
$ws.Range("D118").Value = "SAS"
$ws.Range("E118").Value = 7
$ws.Range("F118").Value = "CC"
$ws.Range("G118").Value = "65464"
$ws.Range("H118").Value = 4031339

# Back to row 111: REFERENCIA filled in last among this batch
$ws.Range("G111").Value = "6465454456"

# Row 119
$ws.Range("A119").Value = 3
$ws.Range("B119").Value = 444441
$ws.Range("C119").Value = "CIMAZ"
$ws.Range("D119").Value = "S.A.S"
$ws.Range("E119").Value = 7
$ws.Range("F119").Value = "CC"
$ws.Range("G119").Value = "45656465"
$ws.Range("H119").Value = 525870

# Row 120
$ws.Range("A120").Value = 3
$ws.Range("B120").Value = 4566546546
$ws.Range("C120").Value = "DIGITALTIC"
$ws.Range("D120").Value = "SAS"
$ws.Range("E120").Value = 51
$ws.Range("F120").Value = "CC"
$ws.Range("G120").Value = "46565464"
$ws.Range("H120").Value = 193970

# View state: window geometry + scroll position + selection
$win = $wb.Windows.Item(1)
$win.Left = -28920
$win.Top = -1815
$win.Width = 29040
$win.Height = 16440

$aw = $excel.ActiveWindow
$aw.ScrollRow = 84
$aw.ScrollColumn = 1

$ws.Range("F108").Select()
